# "Added questions for today" - append 8 new practice-question rows to
# column A, right below the existing list (rows 1-29).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "a-data-structure-question",
    "find-first-non-repeating-character-stream-characters",
    "check-divisibility-binary-stream",
    "select-a-random-number-from-stream-with-o1-space",
    "connect-n-ropes-minimum-cost",
    "minimum-sum-squares-characters-counts-given-string-removing-k-characters",
    "median-of-stream-of-integers-running-integers",
    "lru-cache-implementation"
)

$row = 30
foreach ($v in $values) {
    $ws.Cells.Item($row, 1).Value = $v
    $row++
}

$ws.Range("A30:A37").Select()
